$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values for all source rows before any writes (rotation-safe)
$A3 = $ws.Range("A3").Value2
$B3 = $ws.Range("B3").Value2
$D3 = $ws.Range("D3").Value2
$E3 = $ws.Range("E3").Value2
$F3 = $ws.Range("F3").Value2
$G3 = $ws.Range("G3").Value2
$H3 = $ws.Range("H3").Value2
$Q3 = $ws.Range("Q3").Value2
$R3 = $ws.Range("R3").Value2
$A4 = $ws.Range("A4").Value2
$B4 = $ws.Range("B4").Value2
$D4 = $ws.Range("D4").Value2
$E4 = $ws.Range("E4").Value2
$F4 = $ws.Range("F4").Value2
$G4 = $ws.Range("G4").Value2
$H4 = $ws.Range("H4").Value2
$Q4 = $ws.Range("Q4").Value2
$R4 = $ws.Range("R4").Value2
$A5 = $ws.Range("A5").Value2
$B5 = $ws.Range("B5").Value2
$D5 = $ws.Range("D5").Value2
$E5 = $ws.Range("E5").Value2
$F5 = $ws.Range("F5").Value2
$G5 = $ws.Range("G5").Value2
$H5 = $ws.Range("H5").Value2
$Q5 = $ws.Range("Q5").Value2
$R5 = $ws.Range("R5").Value2
$A8 = $ws.Range("A8").Value2
$B8 = $ws.Range("B8").Value2
$D8 = $ws.Range("D8").Value2
$E8 = $ws.Range("E8").Value2
$F8 = $ws.Range("F8").Value2
$G8 = $ws.Range("G8").Value2
$H8 = $ws.Range("H8").Value2
$Q8 = $ws.Range("Q8").Value2
$R8 = $ws.Range("R8").Value2
$A9 = $ws.Range("A9").Value2
$B9 = $ws.Range("B9").Value2
$D9 = $ws.Range("D9").Value2
$E9 = $ws.Range("E9").Value2
$F9 = $ws.Range("F9").Value2
$G9 = $ws.Range("G9").Value2
$H9 = $ws.Range("H9").Value2
$Q9 = $ws.Range("Q9").Value2
$R9 = $ws.Range("R9").Value2
$A10 = $ws.Range("A10").Value2
$B10 = $ws.Range("B10").Value2
$D10 = $ws.Range("D10").Value2
$E10 = $ws.Range("E10").Value2
$F10 = $ws.Range("F10").Value2
$G10 = $ws.Range("G10").Value2
$H10 = $ws.Range("H10").Value2
$Q10 = $ws.Range("Q10").Value2
$R10 = $ws.Range("R10").Value2
$A11 = $ws.Range("A11").Value2
$B11 = $ws.Range("B11").Value2
$D11 = $ws.Range("D11").Value2
$E11 = $ws.Range("E11").Value2
$F11 = $ws.Range("F11").Value2
$G11 = $ws.Range("G11").Value2
$H11 = $ws.Range("H11").Value2
$Q11 = $ws.Range("Q11").Value2
$R11 = $ws.Range("R11").Value2
$A12 = $ws.Range("A12").Value2
$B12 = $ws.Range("B12").Value2
$D12 = $ws.Range("D12").Value2
$E12 = $ws.Range("E12").Value2
$F12 = $ws.Range("F12").Value2
$G12 = $ws.Range("G12").Value2
$H12 = $ws.Range("H12").Value2
$Q12 = $ws.Range("Q12").Value2
$R12 = $ws.Range("R12").Value2

# Write rotated values into target rows
# Row 3 <- original Row 5
$ws.Range("A3").Value = $A5
$ws.Range("B3").Value = $B5
$ws.Range("D3").Value = $D5
$ws.Range("E3").Value = $E5
$ws.Range("F3").Value = $F5
$ws.Range("G3").Value = $G5
$ws.Range("H3").Value = $H5
$ws.Range("Q3").Value = $Q5
$ws.Range("R3").Value = $R5

# Row 4 <- original Row 3
$ws.Range("A4").Value = $A3
$ws.Range("B4").Value = $B3
$ws.Range("D4").Value = $D3
$ws.Range("E4").Value = $E3
$ws.Range("F4").Value = $F3
$ws.Range("G4").Value = $G3
$ws.Range("H4").Value = $H3
$ws.Range("Q4").Value = $Q3
$ws.Range("R4").Value = $R3

# Row 5 <- original Row 4
$ws.Range("A5").Value = $A4
$ws.Range("B5").Value = $B4
$ws.Range("D5").Value = $D4
$ws.Range("E5").Value = $E4
$ws.Range("F5").Value = $F4
$ws.Range("G5").Value = $G4
$ws.Range("H5").Value = $H4
$ws.Range("Q5").Value = $Q4
$ws.Range("R5").Value = $R4

# Row 8 <- original Row 9
$ws.Range("A8").Value = $A9
$ws.Range("B8").Value = $B9
$ws.Range("D8").Value = $D9
$ws.Range("E8").Value = $E9
$ws.Range("F8").Value = $F9
$ws.Range("G8").Value = $G9
$ws.Range("H8").Value = $H9
$ws.Range("Q8").Value = $Q9
$ws.Range("R8").Value = $R9

# Row 9 <- original Row 11
$ws.Range("A9").Value = $A11
$ws.Range("B9").Value = $B11
$ws.Range("D9").Value = $D11
$ws.Range("E9").Value = $E11
$ws.Range("F9").Value = $F11
$ws.Range("G9").Value = $G11
$ws.Range("H9").Value = $H11
$ws.Range("Q9").Value = $Q11
$ws.Range("R9").Value = $R11

# Row 10 <- original Row 8
$ws.Range("A10").Value = $A8
$ws.Range("B10").Value = $B8
$ws.Range("D10").Value = $D8
$ws.Range("E10").Value = $E8
$ws.Range("F10").Value = $F8
$ws.Range("G10").Value = $G8
$ws.Range("H10").Value = $H8
$ws.Range("Q10").Value = $Q8
$ws.Range("R10").Value = $R8

# Row 11 <- original Row 12
$ws.Range("A11").Value = $A12
$ws.Range("B11").Value = $B12
$ws.Range("D11").Value = $D12
$ws.Range("E11").Value = $E12
$ws.Range("F11").Value = $F12
$ws.Range("G11").Value = $G12
$ws.Range("H11").Value = $H12
$ws.Range("Q11").Value = $Q12
$ws.Range("R11").Value = $R12

# Row 12 <- original Row 10
$ws.Range("A12").Value = $A10
$ws.Range("B12").Value = $B10
$ws.Range("D12").Value = $D10
$ws.Range("E12").Value = $E10
$ws.Range("F12").Value = $F10
$ws.Range("G12").Value = $G10
$ws.Range("H12").Value = $H10
$ws.Range("Q12").Value = $Q10
$ws.Range("R12").Value = $R10

